# Add a new entry to the "Log" sheet (row 43) and to the "Versiones" sheet (row 14)
# describing the favicon L&P feature, per commit v1.12.

$wb = $excel.ActiveWorkbook

# --- Sheet "Log": append new row 43 ---
$logSheet = $wb.Worksheets.Item("Log")

$logSheet.Cells.Item(43, 1).Value = "28/02/2025"
$logSheet.Cells.Item(43, 2).Value = "09:30"
$logSheet.Cells.Item(43, 3).Value = "Favicon L&P en pestaña del navegador"
$logSheet.Cells.Item(43, 4).Value = "Favicon favicon.svg: círculo azul oscuro (#0d2137), texto L&P en blanco, más grande. Enlace en dashboard para que se vea en la solapa del explorador."
$logSheet.Cells.Item(43, 5).Value = "Diagnostico"

# --- Sheet "Versiones": append new row 14 ---
$versionesSheet = $wb.Worksheets.Item("Versiones")

# Leading apostrophe forces the "Versión" value to be stored as text (like the
# other rows' 1.1, 1.10, 1.11, etc.) instead of being auto-converted to the
# number 1.12 by Excel's type inference.
$versionesSheet.Cells.Item(14, 1).Value = "'1.12"
$versionesSheet.Cells.Item(14, 2).Value = "28/02/2025"
$versionesSheet.Cells.Item(14, 3).Value = "Favicon L&P: ícono en pestaña del navegador (fondo azul oscuro, texto blanco)"
